# Apply the "ADC" sheet addition to the workbook.
$wb = $excel.ActiveWorkbook

# Add a new worksheet named "ADC" after the existing "PWM" sheet.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "ADC"

# Row 1: Referenz, V
$ws.Range("A1").Value = "Referenz, V"
$ws.Range("B1").Value = 5

# Row 2: ADC Wert
$ws.Range("A2").Value = "ADC Wert"
$ws.Range("B2").Value = 543
$ws.Range("C2").Value = 566

# Row 3: ADC-Spannung, V
$ws.Range("A3").Value = "ADC-Spannung, V"
$ws.Range("B3").Formula = "=B2*(B1/1024)"

# Row 4: Spannungsteiler
$ws.Range("A4").Value = "Spannungsteiler"
$ws.Range("B4").Value = 4.9000000000000004

# Row 5: echte Spannung, V
$ws.Range("A5").Value = "echte Spannung, V"
$ws.Range("B5").Formula = "=B3*B4"

# Row 4, column C (written after row 5 so the shared-string table order
# matches the target: "echte Spannung, V" precedes "Gesamtwiderstand / Massenwiderstand")
$ws.Range("C4").Value = "Gesamtwiderstand / Massenwiderstand"

# Row 9: echte Spannung, V
$ws.Range("A9").Value = "echte Spannung, V"
$ws.Range("B9").Value = 13

# Row 10: ADC-Spannung, V
$ws.Range("A10").Value = "ADC-Spannung, V"
$ws.Range("B10").Formula = "=B9/B4"

# Row 11: ADC-Wert
$ws.Range("A11").Value = "ADC-Wert"
$ws.Range("B11").Formula = "=(B10*1024)/5"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 17.5703125
$ws.Columns.Item(3).ColumnWidth = 36.85546875

# Page margins: top/bottom = 2 cm (converted to points), left/right stay default (0.7in)
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995

# Select B3 as active cell (matches the saved selection on the ADC sheet)
$ws.Range("B3").Select()

# Restore "PWM" as the active/selected sheet so the workbook re-opens on it,
# same as before the edit.
$pwm = $wb.Worksheets.Item("PWM")
$pwm.Activate()

# Adjust workbook window vertical position
$excel.ActiveWindow.Top = 1800
